$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 18 (pushes existing rows 18..87 down to 19..88,
# and extends the used range from A1:R87 to A1:R88).
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new weekly price-report entry.
$ws.Cells.Item(18, 1).Value = 9
$ws.Cells.Item(18, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(18, 3).Value = "Metropolitana"
$ws.Cells.Item(18, 4).Value = 44620
$ws.Cells.Item(18, 5).Value = 13
$ws.Cells.Item(18, 6).Value = 100114007
$ws.Cells.Item(18, 7).Value = "Jengibre"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 790
$ws.Cells.Item(18, 11).Value = 15000
$ws.Cells.Item(18, 12).Value = 16000
$ws.Cells.Item(18, 13).Value = 15494
$ws.Cells.Item(18, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(18, 15).Value = "Perú"
$ws.Cells.Item(18, 16).Value = 1192
$ws.Cells.Item(18, 17).Value = 13
$ws.Cells.Item(18, 18).Value = "Hortaliza"
